$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# In the "getAssetsinShow" test data block, swap the two existing shownames
# (row 22 becomes "Buffalo Photos", row 23 becomes "Photoshoot")...
$ws.Range("C22").Value = "Buffalo Photos"
$ws.Range("C23").Value = "Photoshoot"

# ...and insert a new row for a third iteration ("Jasper"), which pushes the
# following "CloseDbConnection" block down by one row.
$ws.Rows.Item(24).Insert()

# Copy the formatting/style used by the sibling data rows (e.g. A22) onto the
# new row's first cell so it matches the existing look (centered, quote-prefixed).
$ws.Range("A22").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's values. The leading apostrophe keeps the iteration
# number stored as text (matching the existing "1"/"2" text values) instead
# of a numeric value.
$ws.Range("A24").Value = "'3"
$ws.Range("B24").Value = "Y"
$ws.Range("C24").Value = "Jasper"

# Update the sheet's active selection to the new location of the
# "CloseDbConnection" section header.
$ws.Range("A25").Select()
